$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "Dataset id" in H1
$ws.Range("H1").Value = "Dataset id"

# Add "Dataset id" values for each data row (matching original rows by label in column A)
$ws.Range("H2").Value = 725
$ws.Range("H3").Value = 725
$ws.Range("H4").Value = 726
$ws.Range("H5").Value = 727
$ws.Range("H6").Value = 729
$ws.Range("H7").Value = 728

# Update the active selection to reflect where the cursor ends up after edits
$ws.Range("H8").Select()
